$d = $word.ActiveDocument

# 1. Update the measurement-of-uncertainty VAF/CV% statistics in the
#    "Test Limitations" paragraph.
$d.Content.Find.Execute(
    "5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4%, respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8%, respectively.",
    2
) | Out-Null

# 2. Update the report "saved" date shown in the header/signature block.
$d.Content.Find.Execute(
    "25-Oct-2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "16-Nov-2023",
    2
) | Out-Null
